$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" sheet by copying "2022-Q3" ---
# (copy keeps identical column layout / header style / borders, then we
#  trim it down to the right size and overwrite the data cells)
$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# drop the extra rows copied from 2022-Q3 (13 data rows -> keep header + 7)
$q4Sheet.Range("A9:H13").Delete()

# force columns B,D,E,F,G to be stored as text (matches source formatting,
# e.g. fund codes with leading zeros and decimal figures kept as text)
$q4Sheet.Range("B2:B8").NumberFormat = "@"
$q4Sheet.Range("D2:G8").NumberFormat = "@"

# --- Step 2: write the 2022-Q4 fund holdings data ---
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "002601"
$q4Sheet.Range("C2").Value = "中银证券价值精选灵活配置混合"
$q4Sheet.Range("D2").Value = "5.15"
$q4Sheet.Range("E2").Value = "92.90"
$q4Sheet.Range("F2").Value = "3.82"
$q4Sheet.Range("G2").Value = "0.1967"
$q4Sheet.Range("H2").Value = 9

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "001541"
$q4Sheet.Range("C3").Value = "汇添富民营新动力股票"
$q4Sheet.Range("D3").Value = "3.63"
$q4Sheet.Range("E3").Value = "81.10"
$q4Sheet.Range("F3").Value = "2.46"
$q4Sheet.Range("G3").Value = "0.0893"
$q4Sheet.Range("H3").Value = 10

$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "050014"
$q4Sheet.Range("C4").Value = "博时创业成长混合A"
$q4Sheet.Range("D4").Value = "1.41"
$q4Sheet.Range("E4").Value = "83.80"
$q4Sheet.Range("F4").Value = "3.38"
$q4Sheet.Range("G4").Value = "0.0477"
$q4Sheet.Range("H4").Value = 6

$q4Sheet.Range("A5").Value = 3
$q4Sheet.Range("B5").Value = "011270"
$q4Sheet.Range("C5").Value = "中银证券优势制造股票C"
$q4Sheet.Range("D5").Value = "1.09"
$q4Sheet.Range("E5").Value = "93.20"
$q4Sheet.Range("F5").Value = "4.19"
$q4Sheet.Range("G5").Value = "0.0457"
$q4Sheet.Range("H5").Value = 6

$q4Sheet.Range("A6").Value = 4
$q4Sheet.Range("B6").Value = "011269"
$q4Sheet.Range("C6").Value = "中银证券优势制造股票A"
$q4Sheet.Range("D6").Value = "0.82"
$q4Sheet.Range("E6").Value = "93.20"
$q4Sheet.Range("F6").Value = "4.19"
$q4Sheet.Range("G6").Value = "0.0344"
$q4Sheet.Range("H6").Value = 6

$q4Sheet.Range("A7").Value = 5
$q4Sheet.Range("B7").Value = "002149"
$q4Sheet.Range("C7").Value = "嘉实新优选灵活配置混合"
$q4Sheet.Range("D7").Value = "0.52"
$q4Sheet.Range("E7").Value = "80.11"
$q4Sheet.Range("F7").Value = "5.08"
$q4Sheet.Range("G7").Value = "0.0264"
$q4Sheet.Range("H7").Value = 8

$q4Sheet.Range("A8").Value = 6
$q4Sheet.Range("B8").Value = "002553"
$q4Sheet.Range("C8").Value = "博时创业成长混合C"
$q4Sheet.Range("D8").Value = "0.07"
$q4Sheet.Range("E8").Value = "83.80"
$q4Sheet.Range("F8").Value = "3.38"
$q4Sheet.Range("G8").Value = "0.0024"
$q4Sheet.Range("H8").Value = 6

# --- Step 3: update the "总计" (summary) sheet ---
# row 2 becomes the new 2022-Q4 entry; rows 3-7 shift down by one quarter
# and a new row 8 is appended for the quarter that used to sit in row 7

# give the newly appended index cell (A8) the same style as the other
# numbering cells in column A (e.g. A7) before filling it in
$totalSheet.Range("A7").Copy()
$totalSheet.Range("A8").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.44

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 12
$totalSheet.Range("D3").Value = 0.32

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 6
$totalSheet.Range("D4").Value = 0.19

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2022-Q1"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.09

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q2"
$totalSheet.Range("C6").Value = 7
$totalSheet.Range("D6").Value = 1.55

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q1"
$totalSheet.Range("C7").Value = 3
$totalSheet.Range("D7").Value = 0.65

$totalSheet.Range("A8").Value = 6
$totalSheet.Range("B8").Value = "2020-Q4"
$totalSheet.Range("C8").Value = 9
$totalSheet.Range("D8").Value = 2.36

# leave the view on the summary sheet, like the source workbook
$totalSheet.Activate()